$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("complete", 3210000000),
    @("perfect", 4650000000),
    @("play", 9160000000),
    @("curriculum", 339000000),
    @("older", 1620000000),
    @("young", 3930000000),
    @("serve", 3220000000),
    @("mridangam", 1050000),
    @("response", 2450000000),
    @("shoot", 649000000),
    @("tesla", 396000000),
    @("simple", 3900000000),
    @("science", 2520000000),
    @("solve", 2030000000),
    @("tell", 2100000000),
    @("young", 4220000000),
    @("perfect", 4940000000),
    @("cultural", 1060000000),
    @("serve", 2870000000),
    @("mridangam", 989000),
    @("aspect", 446000000),
    @("tesla", 333000000)
)

$startRow = 213
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
